# Generate Report for Handoff
# Rewrites the Overview / zh-cn / de-de sheets to report on the
# calleeMd1.md, calleeMd2.md, callerMd1.md, callerMd2.md handoff set
# (replacing the previous 3f9c3bde.../705dbe97.../98a5f9f3... file set).

$wb = $excel.ActiveWorkbook

$repoBase   = "https://github.com/OpenLocalizationTest/oltest/blob/398eadb6ab95a69ea4a375472a48f7f6f41db64a/e2e"
$zhCnBase   = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/afb3c439b5fb47ed8afe7cff3a4bb740150ea9a3/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht"
$deDeBase   = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2566fc82aea234bead984a652610e40a048810ad/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht"

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Hyperlinks.Delete()

$overviewRows = @(
    @{ Row = 2; File = "calleeMd1.md" },
    @{ Row = 3; File = "calleeMd2.md" },
    @{ Row = 4; File = "callerMd1.md" },
    @{ Row = 5; File = "callerMd2.md" }
)

foreach ($r in $overviewRows) {
    $row = $r.Row
    $ws1.Cells.Item($row, 2).Value2 = "Ready for handoff"
    $ws1.Cells.Item($row, 3).Value2 = "Ready for handoff"
    $ws1.Cells.Item($row, 4).NumberFormat = "yyyy-mm-dd HH:mm:ss"
    $ws1.Cells.Item($row, 4).Value2 = "2016-03-21 12:54:33"

    $anchor = $ws1.Cells.Item($row, 1)
    [void]$ws1.Hyperlinks.Add($anchor, "$repoBase/$($r.File)", [Type]::Missing, [Type]::Missing, $r.File)
}

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Hyperlinks.Delete()

$zhRows = @(
    @{ Row = 2; File = "calleeMd1.md"; Xlf = "calleeMd1.e8f5ecec2b522eb147a4ff0ca19ca72e17f2186d.zh-cn.xlf"; Dep = "e2e\callerMd2.md,`ne2e\callerMd1.md";    Ref = $null },
    @{ Row = 3; File = "calleeMd2.md"; Xlf = "calleeMd2.63b76063f058ecc63ff1dda71ea2a67db72ae6e1.zh-cn.xlf"; Dep = "e2e\callerMd1.md"; Ref = $null },
    @{ Row = 4; File = "callerMd1.md"; Xlf = "callerMd1.a3bf9f4e7fa2750ec06df0b78a76ae5cafa0e0fd.zh-cn.xlf"; Dep = $null; Ref = "e2e\calleeMd1.md,`ne2e\calleeMd2.md" },
    @{ Row = 5; File = "callerMd2.md"; Xlf = "callerMd2.c7d976edeb9cd5406eae7aba4c05d6d92e81ae95.zh-cn.xlf"; Dep = $null; Ref = "e2e\calleeMd1.md" }
)

foreach ($r in $zhRows) {
    $row = $r.Row
    $ws2.Cells.Item($row, 2).Value2 = ".md"
    $ws2.Cells.Item($row, 3).Value2 = "Ready for handoff"
    $ws2.Cells.Item($row, 5).NumberFormat = "yyyy-mm-dd HH:mm:ss"
    $ws2.Cells.Item($row, 5).Value2 = "2016-03-21 12:54:30"
    $ws2.Cells.Item($row, 8).NumberFormat = "yyyy-mm-dd HH:mm:ss"
    $ws2.Cells.Item($row, 8).Value2 = "0001-01-01 00:00:00"

    # Columns I (Reference Tokens) and K (Dependency From) are mutually
    # exclusive per row in this report; clear any stale value left over
    # from the previous report contents before (re)populating.
    $ws2.Cells.Item($row, 9).ClearContents()
    $ws2.Cells.Item($row, 11).ClearContents()
    if ($r.Ref) {
        $ws2.Cells.Item($row, 9).Value2 = $r.Ref
    }
    $ws2.Cells.Item($row, 10).Value2 = "Include"
    if ($r.Dep) {
        $ws2.Cells.Item($row, 11).Value2 = $r.Dep
    }

    $anchorA = $ws2.Cells.Item($row, 1)
    [void]$ws2.Hyperlinks.Add($anchorA, "$repoBase/$($r.File)", [Type]::Missing, [Type]::Missing, $r.File)

    $anchorD = $ws2.Cells.Item($row, 4)
    [void]$ws2.Hyperlinks.Add($anchorD, "$zhCnBase/$($r.Xlf)", [Type]::Missing, [Type]::Missing, $r.Xlf)
}

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Hyperlinks.Delete()

$deRows = @(
    @{ Row = 2; File = "calleeMd1.md"; Xlf = "calleeMd1.e8f5ecec2b522eb147a4ff0ca19ca72e17f2186d.de-de.xlf"; Dep = "e2e\callerMd2.md,`ne2e\callerMd1.md";    Ref = $null },
    @{ Row = 3; File = "calleeMd2.md"; Xlf = "calleeMd2.63b76063f058ecc63ff1dda71ea2a67db72ae6e1.de-de.xlf"; Dep = "e2e\callerMd1.md"; Ref = $null },
    @{ Row = 4; File = "callerMd1.md"; Xlf = "callerMd1.a3bf9f4e7fa2750ec06df0b78a76ae5cafa0e0fd.de-de.xlf"; Dep = $null; Ref = "e2e\calleeMd1.md,`ne2e\calleeMd2.md" },
    @{ Row = 5; File = "callerMd2.md"; Xlf = "callerMd2.c7d976edeb9cd5406eae7aba4c05d6d92e81ae95.de-de.xlf"; Dep = $null; Ref = "e2e\calleeMd1.md" }
)

foreach ($r in $deRows) {
    $row = $r.Row
    $ws3.Cells.Item($row, 2).Value2 = ".md"
    $ws3.Cells.Item($row, 3).Value2 = "Ready for handoff"
    $ws3.Cells.Item($row, 5).NumberFormat = "yyyy-mm-dd HH:mm:ss"
    $ws3.Cells.Item($row, 5).Value2 = "2016-03-21 12:54:33"
    $ws3.Cells.Item($row, 8).NumberFormat = "yyyy-mm-dd HH:mm:ss"
    $ws3.Cells.Item($row, 8).Value2 = "0001-01-01 00:00:00"

    # Columns I (Reference Tokens) and K (Dependency From) are mutually
    # exclusive per row in this report; clear any stale value left over
    # from the previous report contents before (re)populating.
    $ws3.Cells.Item($row, 9).ClearContents()
    $ws3.Cells.Item($row, 11).ClearContents()
    if ($r.Ref) {
        $ws3.Cells.Item($row, 9).Value2 = $r.Ref
    }
    $ws3.Cells.Item($row, 10).Value2 = "Include"
    if ($r.Dep) {
        $ws3.Cells.Item($row, 11).Value2 = $r.Dep
    }

    $anchorA = $ws3.Cells.Item($row, 1)
    [void]$ws3.Hyperlinks.Add($anchorA, "$repoBase/$($r.File)", [Type]::Missing, [Type]::Missing, $r.File)

    $anchorD = $ws3.Cells.Item($row, 4)
    [void]$ws3.Hyperlinks.Add($anchorD, "$deDeBase/$($r.Xlf)", [Type]::Missing, [Type]::Missing, $r.Xlf)
}

Write-Host "Report regenerated for calleeMd1/calleeMd2/callerMd1/callerMd2"
